$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 5 data rows immediately below the header
# (004526450/MSD, 004431546/GABRIELA, 004448303/NASSIM, 004395314/MARIA, 004240014/ISABELE)
$ws.Range("A2:A6").EntireRow.Delete()

# Insert a new row right after the "005152037 / RODRIGO / 562.4" row
# (that row is now at row 16 following the deletion above) and populate it
$ws.Range("A17").EntireRow.Insert()

$accountCell = $ws.Cells.Item(17, 1)
$accountCell.NumberFormat = "@"
$accountCell.Value = "008002502"
$ws.Cells.Item(17, 2).Value = "JORGEANA"
$ws.Cells.Item(17, 3).Value = 550
